$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells keep their original text formatting (prices stored as text,
# e.g. "1.00" or "331.00", so Excel does not silently coerce them to numbers).
$dCells = @("D2","D3","D4","D5","D6","D8","D9","D10","D11","D12","D15","D16","D17","D18","D19","D20","D22","D23","D24","D28","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D45","D46","D47","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "51.209.24"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "2.728.76"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "115.65"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "331.00"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "41.24"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "20.37"
$ws.Range("D12").Value = "0.0827"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").Value = "3.153.94"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").Value = "2.720.41"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "0.884"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "51.123.79"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").Value = "13.81"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").Value = "3.00"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "288.51"
$ws.Range("E23").Value = "  +3.58%  "
$ws.Range("D24").Value = "70.30"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "10.33"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.141"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "35.84"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "50.05"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "5.61"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").Value = "0.0828"
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("D35").Value = "19.51"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "5.04"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "2.10"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "3.24"
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").Value = "23.87"
$ws.Range("E40").Value = "  +7.52%  "
$ws.Range("D41").Value = "128.87"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("E42").Value = "  +9.07%  "
$ws.Range("E43").Value = "  +3.80%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Value = "3.43"
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("D46").Value = "2.114.07"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Value = "2.19"
$ws.Range("E47").Value = "  +9.90%  "
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("D50").Value = "9.06"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").Value = "60.27"
$ws.Range("E51").Value = "  +0.73%  "
